{"js": "// Applies the 531.DOCX protocol-fill-in edit:\n//  - fills in the date, insulation type, and the first three measurement\n//    rows (u/i/p/R for hh1..hh3) with their real values\n//  - blanks out the still-unfilled placeholder rows (hh4..hh10)\n//\n// Each placeholder token is unique in the document body, so a plain\n// case-sensitive search-and-replace (via body.search) is unambiguous.\n\nasync function replaceOnce(context, needle, replacement) {\n  const results = context.document.body.search(needle, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Placeholder not found: ${needle}`);\n  }\n  // Tokens are unique, but guard against unexpected duplicates anyway.\n  results.items.forEach((range) => {\n    range.insertText(replacement, Word.InsertLocation.replace);\n  });\n  await context.sync();\n}\n\n// Fields that receive real values.\nconst fillIns = [\n  [\"date\", \"16.11.2022\"],\n  [\"tipdv\", \"\u0441\u0442\u0430\u043d\u0434\u0430\u0440\u0442\"],\n  [\"24,3299999237061\", \"24,33\"],\n\n  [\"u1hh\", \"494,9\"],\n  [\"i1hh\", \"6,99\"],\n  [\"p1hh\", \"1665,3\"],\n  [\"R1hh\", \"10\"],\n\n  [\"u2hh\", \"494,1\"],\n  [\"i2hh\", \"7,01\"],\n  [\"p2hh\", \"1644,2\"],\n  [\"R2hh\", \"20\"],\n\n  [\"u3hh\", \"493,7\"],\n  [\"i3hh\", \"7\"],\n  [\"p3hh\", \"1645,7\"],\n  [\"R3hh\", \"30\"],\n];\n\n// Still-unmeasured rows (4..10): the placeholder run is cleared out\n// entirely, leaving an empty paragraph in the cell.\nconst blanks = [];\nfor (let n = 4; n <= 10; n++) {\n  for (const p of [\"u\", \"i\", \"p\", \"R\"]) {\n    blanks.push(`${p}${n}hh`);\n  }\n}\n\nfor (const [needle, replacement] of fillIns) {\n  await replaceOnce(context, needle, replacement);\n}\n\nfor (const needle of blanks) {\n  await replaceOnce(context, needle, \"\");\n}\n", "ps1": "# Applies the 531.DOCX protocol-fill-in edit:\n#  - fills in the date, insulation type, and the first three measurement\n#    rows (u/i/p/R for hh1..hh3) with their real values\n#  - blanks out the still-unfilled placeholder rows (hh4..hh10)\n#\n# Each placeholder token is unique in the document body, so Find/Replace\n# (wdReplaceOne) on the whole-document range is unambiguous.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Placeholder($needle, $replacement) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $needle\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replacement\n    $found = $find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)\n    if (-not $found) {\n        throw \"Placeholder not found: $needle\"\n    }\n}\n\n# Fields that receive real values.\n$fillIns = @(\n    @(\"date\", \"16.11.2022\"),\n    @(\"tipdv\", \"\u0441\u0442\u0430\u043d\u0434\u0430\u0440\u0442\"),\n    @(\"24,3299999237061\", \"24,33\"),\n\n    @(\"u1hh\", \"494,9\"),\n    @(\"i1hh\", \"6,99\"),\n    @(\"p1hh\", \"1665,3\"),\n    @(\"R1hh\", \"10\"),\n\n    @(\"u2hh\", \"494,1\"),\n    @(\"i2hh\", \"7,01\"),\n    @(\"p2hh\", \"1644,2\"),\n    @(\"R2hh\", \"20\"),\n\n    @(\"u3hh\", \"493,7\"),\n    @(\"i3hh\", \"7\"),\n    @(\"p3hh\", \"1645,7\"),\n    @(\"R3hh\", \"30\")\n)\n\nforeach ($pair in $fillIns) {\n    Replace-Placeholder $pair[0] $pair[1]\n}\n\n# Still-unmeasured rows (4..10): the placeholder run is cleared out\n# entirely, leaving an empty paragraph in the cell.\nfor ($n = 4; $n -le 10; $n++) {\n    foreach ($p in @(\"u\", \"i\", \"p\", \"R\")) {\n        $token = \"$p${n}hh\"\n        Replace-Placeholder $token \"\"\n    }\n}\n"}
